# Weekly update: insert this week's two new price rows (Pintón / Primera Pintón)
# at the top of the "Terminal Hortofrutícola Agro Chillán - Plátano" date block,
# pushing the previously-existing rows (349-410) down to (351-412).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 349; everything below shifts down by 2.
$ws.Rows("349:350").Insert()

# --- Row 349: new "Pintón" record for this week ---
$ws.Cells.Item(349, 1).Value = 7
$ws.Cells.Item(349, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(349, 3).Value = "Ñuble"
$ws.Cells.Item(349, 4).Value = 44522
$ws.Cells.Item(349, 5).Value = 16
$ws.Cells.Item(349, 6).Value = "Fruta"
$ws.Cells.Item(349, 7).Value = 100108
$ws.Cells.Item(349, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(349, 9).Value = 100108006
$ws.Cells.Item(349, 10).Value = "Plátano"
$ws.Cells.Item(349, 11).Value = "Sin especificar"
$ws.Cells.Item(349, 12).Value = "Pintón"
$ws.Cells.Item(349, 13).Value = 60
$ws.Cells.Item(349, 14).Value = 15000
$ws.Cells.Item(349, 15).Value = 15000
$ws.Cells.Item(349, 16).Value = 15000
$ws.Cells.Item(349, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(349, 18).Value = "Ecuador"
$ws.Cells.Item(349, 19).Value = 750
$ws.Cells.Item(349, 20).Value = 20

# --- Row 350: new "Primera Pintón" record for this week ---
$ws.Cells.Item(350, 1).Value = 7
$ws.Cells.Item(350, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(350, 3).Value = "Ñuble"
$ws.Cells.Item(350, 4).Value = 44522
$ws.Cells.Item(350, 5).Value = 16
$ws.Cells.Item(350, 6).Value = "Fruta"
$ws.Cells.Item(350, 7).Value = 100108
$ws.Cells.Item(350, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(350, 9).Value = 100108006
$ws.Cells.Item(350, 10).Value = "Plátano"
$ws.Cells.Item(350, 11).Value = "Sin especificar"
$ws.Cells.Item(350, 12).Value = "Primera Pintón"
$ws.Cells.Item(350, 13).Value = 160
$ws.Cells.Item(350, 14).Value = 16000
$ws.Cells.Item(350, 15).Value = 17000
$ws.Cells.Item(350, 16).Value = 16500
$ws.Cells.Item(350, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(350, 18).Value = "Ecuador"
$ws.Cells.Item(350, 19).Value = 825
$ws.Cells.Item(350, 20).Value = 20

# Make sure the date cells keep the date/time number format used by the rest
# of column D.
$ws.Range("D349:D350").NumberFormat = $ws.Range("D351").NumberFormat
